# Preliminary check-in: Update to some ODK Survey forms to rename them to the
# table_id so that we generate definitions.csv and properties.csv ; update to
# process a properties sheet into the properties.csv ; minimize the content
# of the properties.csv

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "settings" sheet: the setting that used to be named "form_id" is renamed
#    to "table_id" (its value, "scan_example", stays the same). Move the
#    selection to A3 (it will get re-activated below when we add/activate the
#    new "properties" sheet, so the tabSelected flag ends up on that one).
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("A2").Value = "table_id"
[void]$wsSettings.Range("A3").Select()

# ---------------------------------------------------------------------------
# 2) Add a new "properties" sheet after the last existing sheet and populate
#    it with the partition/aspect/key/type/value rows describing the table.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsProps = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsProps.Name = "properties"

$colOrderJson = '["address","address_image0_contentType","address_image0_uriFragment","comments","comments_image0_contentType","comments_image0_uriFragment","fri_chores","fri_chores_image0_contentType","fri_chores_image0_uriFragment","mon_chores","mon_chores_image0_contentType","mon_chores_image0_uriFragment","name","name_image0_contentType","name_image0_uriFragment","qrcode","qrcode_image0_contentType","qrcode_image0_uriFragment","roomNum","roomNum_image0_contentType","roomNum_image0_uriFragment","sat_chores","sat_chores_image0_contentType","sat_chores_image0_uriFragment","scan_output_directory","stay","stay_image0_contentType","stay_image0_uriFragment","sun_chores","sun_chores_image0_contentType","sun_chores_image0_uriFragment","thurs_chores","thurs_chores_image0_contentType","thurs_chores_image0_uriFragment","tues_chores","tues_chores_image0_contentType","tues_chores_image0_uriFragment","wed_chores","wed_chores_image0_contentType","wed_chores_image0_uriFragment"]'

$propRows = @(
    @("partition", "aspect", "key", "type", "value"),
    @("Table", "default", "colOrder", "array", $colOrderJson),
    @("Table", "default", "defaultViewType", "string", "LIST"),
    @("Table", "default", "listViewFileName", "configpath", "config/tables/scan_example/html/scan_example_list.html"),
    @("Table", "default", "detailViewFileName", "configpath", "config/tables/scan_example/html/scan_example_detail.html")
)

for ($r = 0; $r -lt $propRows.Length; $r++) {
    $row = $propRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsProps.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ColumnWidth (character units) round-trips to the stored XML "width" with a
# constant ~0.8333 (5px/6px) padding offset added back in by the engine, so
# back that offset out here to land on the exact target widths (15 / 23).
$padOffset = 0.8333333333333334
$wsProps.Columns.Item(1).ColumnWidth = 15 - $padOffset
$wsProps.Columns.Item(2).ColumnWidth = 15 - $padOffset
$wsProps.Columns.Item(3).ColumnWidth = 23 - $padOffset
$wsProps.Columns.Item(4).ColumnWidth = 15 - $padOffset
$wsProps.Columns.Item(5).ColumnWidth = 15 - $padOffset

# Make "properties" the active/selected sheet, with F10 selected, matching
# the tail state of the edited workbook.
[void]$wsProps.Activate()
[void]$wsProps.Range("F10").Select()
